# Weekly update: insert a new price record at the top of the Zanahoria
# data block (row 262), pushing all existing records down by one row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row before the current row 262 - this shifts rows
# 262:309 down to 263:310 and extends the used range to R310.
$ws.Rows("262:262").Insert()

# Populate the newly inserted row 262 with the new weekly record.
$ws.Range("A262").Value = 7
$ws.Range("B262").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C262").Value = "Ñuble"
$ws.Range("D262").Value = 44776
$ws.Range("E262").Value = 16
$ws.Range("F262").Value = 100114013
$ws.Range("G262").Value = "Zanahoria"
$ws.Range("H262").Value = "Sin especificar"
$ws.Range("I262").Value = "Primera"
$ws.Range("J262").Value = 100
$ws.Range("K262").Value = 8000
$ws.Range("L262").Value = 8500
$ws.Range("M262").Value = 8250
$ws.Range("N262").Value = "$/saco 20 kilos"
$ws.Range("O262").Value = "Provincia de Diguillín"
$ws.Range("P262").Value = 412
$ws.Range("Q262").Value = 20
$ws.Range("R262").Value = "Hortaliza"
